# Auto-generated edit script: updates crypto price/volume columns (D, E)
# to match the target snapshot described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Temporarily force Text format over the data range so that numeric-looking
# strings (e.g. "23.98", "9.16") are preserved as text, matching the
# original inline-string cell types instead of being parsed as numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.589.65"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.576.47"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").Value = "213.26"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "44.91"
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("D9").Value = "23.98"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "0.247"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "1.802.80"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "1.562.70"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").Value = "28.599.23"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "3.69"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "62.25"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").Value = "230.58"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "0.0₃0691"
$ws.Range("E21").Value = "  -2.24%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  -4.63%  "
$ws.Range("D24").Value = "9.16"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").Value = "2.08"
$ws.Range("E25").Value = "  +6.75%  "
$ws.Range("D26").Value = "151.49"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "15.02"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").Value = "6.44"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").Value = "0.0484"
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").Value = "3.10"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D35").Value = "1.400.65"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("E37").Value = "  -3.69%  "
$ws.Range("D38").Value = "2.36"
$ws.Range("E39").Value = "  +2.79%  "
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").Value = "0.523"
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").Value = "0.792"
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("D44").Value = "1.88"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "0.0466"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("D48").Value = "62.92"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").Value = "1.714.63"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "86.48"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  -0.44%  "

# Restore the original (default/Normal) style so no stray style index is
# introduced on the touched cells.
$dataRange.Style = "Normal"

Write-Host "Applied 78 cell updates to columns D and E (rows 2-51)."
